$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.901.74'
$ws.Range('E2').Value = '  +2.80%  '
$ws.Range('D3').Value = '1.713.49'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'311.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = "'0.3776"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').Value = "'49.56"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.23%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = "'0.3484"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').Value = "'0.07498"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.08%  '
$ws.Range('D12').Value = "'1.004"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = "'20.95"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('D14').Value = "'6.296"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.43%  '
$ws.Range('E15').Value = '  +3.70%  '
$ws.Range('D16').Value = '1.713.89'
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = "'0.00001127"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = "'1.001"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = "'0.06733"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = "'84.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.56%  '
$ws.Range('D21').Value = "'17.26"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.17%  '
$ws.Range('D22').Value = "'6.387"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.22%  '
$ws.Range('D23').Value = "'13.11"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.03%  '
$ws.Range('D24').Value = '24.826.67'
$ws.Range('E24').Value = '  +2.64%  '
$ws.Range('D25').Value = "'2.444"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').Value = "'2.793"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.26%  '
$ws.Range('D27').Value = "'20.49"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.30%  '
$ws.Range('D28').Value = "'150.93"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').Value = "'132.23"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.21%  '
$ws.Range('D30').Value = '1.908.75'
$ws.Range('E30').Value = '  +2.55%  '
$ws.Range('D31').Value = "'1.177"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +19.30%  '
$ws.Range('D32').Value = "'6.842"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.55%  '
$ws.Range('D33').Value = "'4.237"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.39%  '
$ws.Range('D34').Value = "'13.84"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.19%  '
$ws.Range('D35').Value = "'0.08834"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.57%  '
$ws.Range('D36').Value = "'1.769"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('D38').Value = "'0.06563"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('D39').Value = "'0.02403"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.82%  '
$ws.Range('D40').Value = "'8.997"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('D41').Value = "'0.2201"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.42%  '
$ws.Range('D42').Value = "'1.280"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').Value = "'0.6431"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.46%  '
$ws.Range('D45').Value = "'13.94"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.93%  '
$ws.Range('D46').Value = "'0.6142"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.39%  '
$ws.Range('D47').Value = "'3.822"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').Value = "'2.138"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.90%  '
$ws.Range('D49').Value = "'130.13"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.60%  '
$ws.Range('D50').Value = "'0.07281"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('D51').Value = "'79.90"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.38%  '
